# Auto-generated edit script applying cryptos.xlsx diff (Mon May 13 07:35:00 UTC 2024 update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.667.09'
$ws.Range("E2").Value = '  +1.44%  '

$ws.Range("D3").Value = '2.926.26'
$ws.Range("E3").Value = '  +0.59%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''595.55'
$ws.Range("E5").Value = '  +1.05%  '

$ws.Range("D6").Value = '''141.58'
$ws.Range("E6").Value = '  -1.91%  '

$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '2.924.30'
$ws.Range("E8").Value = '  +0.55%  '

$ws.Range("E9").Value = '  -1.23%  '

$ws.Range("D10").Value = '''7.18'
$ws.Range("E10").Value = '  +4.59%  '

$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("D12").Value = '''0.441'
$ws.Range("E12").Value = '  +0.55%  '

$ws.Range("E13").Value = '  -0.97%  '

$ws.Range("D14").Value = '''33.02'
$ws.Range("E14").Value = '  -0.99%  '

$ws.Range("E15").Value = '  -0.09%  '

$ws.Range("D16").Value = '3.411.46'
$ws.Range("E16").Value = '  +0.55%  '

$ws.Range("D17").Value = '61.552.34'
$ws.Range("E17").Value = '  +1.31%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '''6.64'
$ws.Range("E18").Value = '  -0.35%  '

$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '2.916.93'
$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").Value = '''433.25'
$ws.Range("E20").Value = '  +0.71%  '

$ws.Range("D21").Value = '''13.46'
$ws.Range("E21").Value = '  +0.92%  '

$ws.Range("D22").Value = '''0.670'
$ws.Range("E22").Value = '  -0.80%  '

$ws.Range("D23").Value = '''7.04'
$ws.Range("E23").Value = '  -0.64%  '

$ws.Range("D24").Value = '''81.02'
$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("D25").Value = '''10.64'
$ws.Range("E25").Value = '  -2.25%  '

$ws.Range("D26").Value = '''2.11'
$ws.Range("E26").Value = '  -3.17%  '

$ws.Range("D27").Value = '''11.68'
$ws.Range("E27").Value = '  -0.55%  '

$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("E29").Value = '  -0.14%  '

$ws.Range("E30").Value = '  -7.84%  '

$ws.Range("E31").Value = '  -1.57%  '

$ws.Range("D32").Value = '''26.15'
$ws.Range("E32").Value = '  -1.38%  '

$ws.Range("E33").Value = '  +0.06%  '

$ws.Range("D34").Value = '''0.105'
$ws.Range("E34").Value = '  -2.94%  '

$ws.Range("D35").Value = '0.0₃0855'
$ws.Range("E35").Value = '  -0.60%  '

$ws.Range("D36").Value = '''0.985'
$ws.Range("E36").Value = '  -2.13%  '

$ws.Range("E37").Value = '  -0.87%  '

$ws.Range("E38").Value = '  -0.83%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '''1.94'
$ws.Range("E39").Value = '  -1.43%  '

$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").Value = '''2.83'
$ws.Range("E40").Value = '  -4.93%  '

$ws.Range("D41").Value = '''8.48'

$ws.Range("E42").Value = '  -2.86%  '

$ws.Range("E43").Value = '  -3.09%  '

$ws.Range("D44").Value = '''38.54'
$ws.Range("E44").Value = '  -7.01%  '

$ws.Range("D45").Value = '2.683.39'
$ws.Range("E45").Value = '  -0.38%  '

$ws.Range("D46").Value = '''133.44'
$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("D47").Value = '''0.0336'
$ws.Range("E47").Value = '  -1.85%  '

$ws.Range("D48").Value = '''356.78'
$ws.Range("E48").Value = '  -4.73%  '

$ws.Range("D50").Value = '''22.71'

$ws.Range("E51").Value = '  -2.13%  '
